$d = $word.ActiveDocument

function Get-ParaIndex($searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { return -1 }
    $pre = $d.Range(0, $rng.Start)
    return $pre.Paragraphs.Count + 1
}

# --- Merge split runs into single runs (text content unchanged) ---
$null = $d.Content.Find.Execute("az acr create --resource-group posiorg --name posioapp --sku Basic", $true, $false, $false, $false, $false, $true, 1, $false, "az acr create --resource-group posiorg --name posioapp --sku Basic", 2)
$null = $d.Content.Find.Execute("az acr login --name posioapp ", $true, $false, $false, $false, $false, $true, 1, $false, "az acr login --name posioapp ", 2)
$null = $d.Content.Find.Execute("docker context create aci posiocontext", $true, $false, $false, $false, $false, $true, 1, $false, "docker context create aci posiocontext", 2)
$null = $d.Content.Find.Execute("docker context use posiocontext", $true, $false, $false, $false, $false, $true, 1, $false, "docker context use posiocontext", 2)
$null = $d.Content.Find.Execute("az group delete --name posiorg", $true, $false, $false, $false, $false, $true, 1, $false, "az group delete --name posiorg", 2)

# --- Insert a new blank paragraph right after the "storage_account_name: taskboardstorageacc" line ---
$null = $d.Content.Find.Execute("            storage_account_name: taskboardstorageacc", $true, $false, $false, $false, $false, $true, 1, $false, "            storage_account_name: taskboardstorageacc`r", 2)

# --- Remove one of the two blank paragraphs that follow "az group delete --name taskBoardResourceGroup" ---
$idxC = Get-ParaIndex("az group delete --name taskBoardResourceGroup")
$pC = $d.Paragraphs.Item($idxC)
$pC.Next().Range.Delete()

# --- Remove two of the four blank paragraphs that follow the "az container create ...trackerapprggeorgi" line ---
$idxB = Get-ParaIndex("az container create --resource-group trackerapprggeorgi")
$pB = $d.Paragraphs.Item($idxB)
$n1 = $pB.Next()
$n2 = $n1.Next()
$n2.Range.Delete()
$n1.Range.Delete()

# --- Remove one of the two blank paragraphs that precede the "az container create ...trackerapprggeorgi" line ---
$idxA = Get-ParaIndex("az container create --resource-group trackerapprggeorgi")
$pA = $d.Paragraphs.Item($idxA)
$pA.Previous().Range.Delete()

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
